{"js": "// Rewrite the report body (keeping the title heading paragraph intact) to\n// match the restructured \"System Overview / Patch Status Summary /\n// Compliance with RMF Controls / Recommended next steps / Risk Assessment\"\n// layout described by the diff: bullets get \"* \" / \"1.\" / \"2.\" / \"3.\"\n// markers, blank-line paragraphs are inserted between sections, and\n// several paragraphs are reworded or resplit.\n\nconst body = context.document.body;\n\n// [text, hasLineBreak] for every paragraph that follows the title, in the\n// final order. A trailing \"\\u000b\" (Word's manual line-break char)\n// reproduces the <w:br/> that trails almost every run in the source\n// document; the very last paragraph has no trailing break.\nconst target = [\n  [\"*** System Overview ***\", true],\n  [\"\", true],\n  [\"The system is a computer with the following specifications:\", true],\n  [\"\", true],\n  [\"* Date: April 6th, 2025\", true],\n  [\"* Time: 4:02 PM\", true],\n  [\"* OS Name: kb322-18\", true],\n  [\"* OS Version: Debian 6.1.129-1 (2025-03-06)\", true],\n  [\"* Computer Name: kb322-18\", true],\n  [\"* IP Address: 140.160.138.147\", true],\n  [\"\", true],\n  [\"*** Patch Status Summary***\", true],\n  [\"\", true],\n  [\"The following patches are pending updates:\", true],\n  [\"\", true],\n  [\"1. Code/stable 1.99.0-1743632463 amd64 [upgradable from: 1.98.2-1741788907]\", true],\n  [\"\", true],\n  [\"There is no information available about the relevance of this patch to security.\", true],\n  [\"\", true],\n  [\"*** Compliance with RMF Controls***\", true],\n  [\"\", true],\n  [\"To ensure compliance with the Risk Management Framework (RMF), we recommend:\", true],\n  [\"\", true],\n  [\"* Flaw remediation in place by installing the pending update as soon as possible.\", true],\n  [\"* Identification, reporting, and corrective action to monitor and report any vulnerabilities.\", true],\n  [\"* Configuration management to track changes to system configurations and ensure consistency.\", true],\n  [\"* Vulnerability checks to regularly scan for new vulnerabilities.\", true],\n  [\"\", true],\n  [\"*** Recommended next steps***\", true],\n  [\"\", true],\n  [\"The recommended next steps are:\", true],\n  [\"\", true],\n  [\"1. Review and assess updates to determine the relevance of the pending patch to security.\", true],\n  [\"2. Scheduling patch deployments to install the update as soon as possible.\", true],\n  [\"3. Guidance for Update documentation to track changes to system configurations and ensure consistency.\", true],\n  [\"\", true],\n  [\"*** Risk Assessment***\", true],\n  [\"\", true],\n  [\"There is currently no information available about potential vulnerabilities, but the pending patch may address a known vulnerability. If installed promptly, this should mitigate any risk associated with the vulnerability. The potential impact level of this vulnerability is unknown, and further assessment would be required after installation of the patch.\", false]\n];\n\n// Paragraph 0 is the \"Operating System Patch Management RMF Compliance\"\n// Heading 2 title -- the diff leaves it untouched. Paragraph 1 (originally\n// \"*** System Overview ***\") is a plain body paragraph with no explicit\n// style, so reuse it as the anchor for the rebuilt content instead of\n// inserting next to the Heading 2 paragraph (which would make new\n// paragraphs inherit the Heading 2 style).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = paragraphs.items.length - 1; i >= 2; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nconst anchorParas = body.paragraphs;\nanchorParas.load(\"items\");\nawait context.sync();\nlet anchor = anchorParas.items[1];\n\nconst [firstText, firstHasBreak] = target[0];\nanchor.insertText(firstHasBreak ? firstText + \"\\u000b\" : firstText, Word.InsertLocation.replace);\nawait context.sync();\n\nlet insertAfter = anchor;\nfor (let i = 1; i < target.length; i++) {\n  const [text, hasBreak] = target[i];\n  const full = hasBreak ? text + \"\\u000b\" : text;\n  insertAfter = insertAfter.insertParagraph(full, Word.InsertLocation.after);\n}\nawait context.sync();\n", "ps1": "# Rewrite the report body (keeping the title Heading 2 paragraph intact) to\n# match the restructured \"System Overview / Patch Status Summary /\n# Compliance with RMF Controls / Recommended next steps / Risk Assessment\"\n# layout described by the diff: bullets get \"* \" / \"1.\" / \"2.\" / \"3.\"\n# markers, blank-line paragraphs are inserted between sections, and\n# several paragraphs are reworded or resplit.\n\n$d = $word.ActiveDocument\n\n# Every paragraph text that follows the title, in final order. [char]11 is\n# Word's manual line-break char (0x0B) -- appending it reproduces the\n# <w:br/> that trails almost every run in the source document; the very\n# last paragraph has no trailing break.\n$target = @(\n  \"*** System Overview ***\" + [char]11,\n  \"\" + [char]11,\n  \"The system is a computer with the following specifications:\" + [char]11,\n  \"\" + [char]11,\n  \"* Date: April 6th, 2025\" + [char]11,\n  \"* Time: 4:02 PM\" + [char]11,\n  \"* OS Name: kb322-18\" + [char]11,\n  \"* OS Version: Debian 6.1.129-1 (2025-03-06)\" + [char]11,\n  \"* Computer Name: kb322-18\" + [char]11,\n  \"* IP Address: 140.160.138.147\" + [char]11,\n  \"\" + [char]11,\n  \"*** Patch Status Summary***\" + [char]11,\n  \"\" + [char]11,\n  \"The following patches are pending updates:\" + [char]11,\n  \"\" + [char]11,\n  \"1. Code/stable 1.99.0-1743632463 amd64 [upgradable from: 1.98.2-1741788907]\" + [char]11,\n  \"\" + [char]11,\n  \"There is no information available about the relevance of this patch to security.\" + [char]11,\n  \"\" + [char]11,\n  \"*** Compliance with RMF Controls***\" + [char]11,\n  \"\" + [char]11,\n  \"To ensure compliance with the Risk Management Framework (RMF), we recommend:\" + [char]11,\n  \"\" + [char]11,\n  \"* Flaw remediation in place by installing the pending update as soon as possible.\" + [char]11,\n  \"* Identification, reporting, and corrective action to monitor and report any vulnerabilities.\" + [char]11,\n  \"* Configuration management to track changes to system configurations and ensure consistency.\" + [char]11,\n  \"* Vulnerability checks to regularly scan for new vulnerabilities.\" + [char]11,\n  \"\" + [char]11,\n  \"*** Recommended next steps***\" + [char]11,\n  \"\" + [char]11,\n  \"The recommended next steps are:\" + [char]11,\n  \"\" + [char]11,\n  \"1. Review and assess updates to determine the relevance of the pending patch to security.\" + [char]11,\n  \"2. Scheduling patch deployments to install the update as soon as possible.\" + [char]11,\n  \"3. Guidance for Update documentation to track changes to system configurations and ensure consistency.\" + [char]11,\n  \"\" + [char]11,\n  \"*** Risk Assessment***\" + [char]11,\n  \"\" + [char]11,\n  \"There is currently no information available about potential vulnerabilities, but the pending patch may address a known vulnerability. If installed promptly, this should mitigate any risk associated with the vulnerability. The potential impact level of this vulnerability is unknown, and further assessment would be required after installation of the patch.\"\n)\n\n# Paragraph 1 is the \"Operating System Patch Management RMF Compliance\"\n# Heading 2 title -- the diff leaves it untouched. Paragraph 2 (originally\n# \"*** System Overview ***\") is a plain body paragraph with no explicit\n# style, so reuse it as the anchor for the rebuilt content instead of\n# inserting next to the Heading 2 paragraph (which would make new\n# paragraphs inherit the Heading 2 style).\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 3; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n$anchor = $d.Paragraphs.Item(2)\n$anchor.Range.Text = $target[0]\n\nfor ($i = 1; $i -lt $target.Length; $i++) {\n    $last = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $last.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $newPara.Range.Text = $target[$i]\n}\n"}
